# Generate Report for Handback
# Adds a new handed-back file ("690a98eb-5eed-4f24-9311-4d0094485c6c.md") as
# row 3 on the Overview / zh-cn / de-de sheets, mirroring the existing
# f83d969b... row, and extends the tables to cover the new row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newFile    = "690a98eb-5eed-4f24-9311-4d0094485c6c.md"
$newPath    = "e2e\690a98eb-5eed-4f24-9311-4d0094485c6c.md"
$newStatus  = "Handed back: in sync with en-US"
$genDate    = "2017-02-22 06:52:01"

$zhXlf      = "690a98eb-5eed-4f24-9311-4d0094485c6c.423573a796c0147d78d10216e4df8427a21823d0.zh-cn.xlf"
$zhHoDate   = "2017-02-22 06:51:45"
$zhHbDate   = "2017-02-22 06:52:42"

$deXlf      = "690a98eb-5eed-4f24-9311-4d0094485c6c.423573a796c0147d78d10216e4df8427a21823d0.de-de.xlf"
$deHoDate   = $genDate
$deHbDate   = "2017-02-22 06:53:05"

# ---------------------------------------------------------------------
# Overview sheet - row 3
# ---------------------------------------------------------------------
$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = $newPath
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G3").Value = $genDate
$ws1.Range("G3").NumberFormat = $ws1.Range("G2").NumberFormat

$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/bdb8d09cd6fb8778b5e72cdb511318b10bbcb2ac/$newPath", "", "", $newPath)
$ws1.Range("B3").Font.Underline = 2
$ws1.Range("B3").Font.Color = 15570276

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet - row 3
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = $newFile
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = $newStatus
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = $zhHoDate
$ws2.Range("H3").NumberFormat = $ws2.Range("H2").NumberFormat
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = $newFile
$ws2.Range("K3").Value = $zhXlf
$ws2.Range("L3").Value = $zhHbDate
$ws2.Range("L3").NumberFormat = $ws2.Range("L2").NumberFormat
$ws2.Range("M3").Value = ""
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "True"
$ws2.Range("P3").Value = ""
$ws2.Range("Q3").Value = "False"
$ws2.Range("R3").Value = ""

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/bdb8d09cd6fb8778b5e72cdb511318b10bbcb2ac/e2e/$newFile", "", "", $newFile)
$ws2.Range("A3").Font.Underline = 2
$ws2.Range("A3").Font.Color = 15570276

$ws2.Hyperlinks.Add($ws2.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/423573a796c0147d78d10216e4df8427a21823d0/e2e/$newFile", "", "", $newFile)
$ws2.Range("J3").Font.Underline = 2
$ws2.Range("J3").Font.Color = 15570276

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:R3"))

# ---------------------------------------------------------------------
# de-de sheet - row 3
# ---------------------------------------------------------------------
$ws3.Range("A3").Value = $newFile
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = $newStatus
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = $deHoDate
$ws3.Range("H3").NumberFormat = $ws3.Range("H2").NumberFormat
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = $newFile
$ws3.Range("K3").Value = $deXlf
$ws3.Range("L3").Value = $deHbDate
$ws3.Range("L3").NumberFormat = $ws3.Range("L2").NumberFormat
$ws3.Range("M3").Value = ""
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "True"
$ws3.Range("P3").Value = ""
$ws3.Range("Q3").Value = "False"
$ws3.Range("R3").Value = ""

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/bdb8d09cd6fb8778b5e72cdb511318b10bbcb2ac/e2e/$newFile", "", "", $newFile)
$ws3.Range("A3").Font.Underline = 2
$ws3.Range("A3").Font.Color = 15570276

$ws3.Hyperlinks.Add($ws3.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/423573a796c0147d78d10216e4df8427a21823d0/e2e/$newFile", "", "", $newFile)
$ws3.Range("J3").Font.Underline = 2
$ws3.Range("J3").Font.Color = 15570276

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:R3"))
